$d = $word.ActiveDocument

function Insert-Break {
    param(
        [string]$Find,
        [string]$Replace
    )
    $rng = $d.Content
    $ok = $rng.Find.Execute($Find, $true, $false, $false, $false, $false, $true, 1, $false, $Replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $Find"
    }
}

Insert-Break "Bibliografia:A bibliografia" "Bibliografia:^lA bibliografia"

Insert-Break "no tópico de mentoria: [1] Peddy" "no tópico de mentoria: ^l[1] Peddy"

Insert-Break "Houston: Bullion Books, 2001. [2] Zachary" "Houston: Bullion Books, 2001. ^l[2] Zachary"

Insert-Break "Lisboa: IST Press, 2011. p. 19-27. [3] Mueller" "Lisboa: IST Press, 2011. p. 19-27. ^l[3] Mueller"

Insert-Break "European Journal of Engineering Education, 2004. [4] Kaul" "European Journal of Engineering Education, 2004. ^l[4] Kaul"

Insert-Break "Education, v. 21, p. 14-23,2019. [5] Diretrizes" "Education, v. 21, p. 14-23,2019. ^l[5] Diretrizes"

Write-Host "Done"
